$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.362.54'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '2.643.49'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '597.87'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.84'
$ws.Range('E6').Value = '  +1.21%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.546'
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '2.642.71'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('E10').Value = '  +8.31%  '
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.25'
$ws.Range('E14').Value = '  +2.32%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000191'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').Value = '3.123.61'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '68.336.56'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '2.651.89'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.45'
$ws.Range('E19').Value = '  +1.43%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '364.85'
$ws.Range('E20').Value = '  -2.05%  '
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.41'
$ws.Range('E22').Value = '  +3.95%  '
$ws.Range('E23').Value = '  +1.90%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.06'
$ws.Range('E24').Value = '  +0.68%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.33'
$ws.Range('E25').Value = '  +2.21%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.84'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('D29').Value = '2.775.95'
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '573.51'
$ws.Range('E31').Value = '  -0.71%  '
$ws.Range('E32').Value = '  +4.62%  '
$ws.Range('E33').Value = '  +2.03%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.86'
$ws.Range('E34').Value = '  +0.88%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.131'
$ws.Range('E35').Value = '  +3.74%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +5.36%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '160.40'
$ws.Range('E38').Value = '  +1.15%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '19.41'
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.374'
$ws.Range('E40').Value = '  +1.33%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.90'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.43'
$ws.Range('E42').Value = '  +2.00%  '
$ws.Range('E43').Value = '  +8.19%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.66'
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('E45').Value = '  +3.65%  '
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '40.36'
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '157.14'
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.72'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '21.98'
$ws.Range('E51').Value = '  +0.70%  '
